$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Création des clips de maps dans le dernier frame
# (shift/renumber of "Forêt perdue" map clip labels ; placeholder
# "Forêt perdue" entries in row 6/7 become distinct numbered clips,
# cascading the earlier numbered clips forward)

$ws.Range("B5").Value = "Forêt perdue 4"
$ws.Range("B6").Value = "Forêt perdue 2"
$ws.Range("C6").Value = "Forêt perdue 3"
$ws.Range("C7").Value = "Forêt perdue 5"
$ws.Range("D6").Value = "Forêt perdue 6"
$ws.Range("B7").Value = "Forêt perdue 7"
$ws.Range("D7").Value = "Forêt perdue 8"

$ws.Range("F6").Select()
